$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the spelling of "serorreversion" -> "seroreversion" in the second
# table's header row (row 9). The displayed text for G9 ("Tasa de
# seroreversión Rhat") is unchanged, but we rewrite it too so the shared
# string pool is rebuilt/reordered the same way Excel would on a real edit.
$ws.Range("F9").Value = "Tasa de seroreversion"
$ws.Range("G9").Value = "Tasa de seroreversión Rhat"

# Update the selection to cover the whole second table (A9:G13) with the
# active cell still at A9.
$ws.Range("A9:G13").Select()
